# Apply "Modif url canonique termino" changes to the workbook.
$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update the generation Date value ---
$wsMetadata = $wb.Worksheets.Item("Metadata")
$wsMetadata.Range("B8").Value = "2025-07-25T07:22:51+00:00"

# --- Elements sheet: update canonical terminology URLs and column width ---
$wsElements = $wb.Worksheets.Item("Elements")

$wsElements.Range("Z3").Value = "https://mos.esante.gouv.fr/NOS/TRE_R14-TypeDiplome/FHIR/TRE-R14-TypeDiplome?vs"
$wsElements.Range("Z4").Value = "https://mos.esante.gouv.fr/NOS/TRE_R16-LieuFormation/FHIR/TRE-R16-LieuFormation?vs"
$wsElements.Range("Z7").Value = "https://mos.esante.gouv.fr/NOS/TRE_R49-DiplomeEtudeSpecialisee/FHIR/TRE-R49-DiplomeEtudeSpecialisee?vs"

# Widen column Z (Binding Value Set) to fit the longer URLs.
$wsElements.Columns.Item(26).ColumnWidth = 88.1953125
